# Update data in excel files
# - "critical" sheet (sheet1): replace rows 2-13 with the 9 rows of data that
#   now also appear on the "regular" sheet, add client_id (B) + report_date (D)
#   columns with matching styles, delete the now-unused trailing rows 11-13.
# - "regular" sheet (sheet2): re-order/update the 9 data rows.

$wb = $excel.ActiveWorkbook
$critical = $wb.Worksheets.Item("critical")
$regular  = $wb.Worksheets.Item("regular")

# ---------------------------------------------------------------------------
# Final data for both sheets, rows 2..10: problem_id, client_id, product_id,
# report_date (serial), description
# ---------------------------------------------------------------------------
$rows = @(
    @(12347, 32145, "3434T", 43789, "Not working"),
    @(12365, 12345, "3625N", 43790, "Not working"),
    @(12346, 96325, "2333B", 43793, "Not working"),
    @(12348, 12347, "3231R", 43794, "Not working"),
    @(12398, 32146, "1236G", 43794, "Not working"),
    @(12399, 12345, "1425F", 43794, "Not working"),
    @(12333, 32145, "1596D", 43794, "Not working"),
    @(12345, 12345, "1234A", 43795, "Not working"),
    @(12121, 96325, "6669R", 43795, "Not working")
)

# ---------------------------------------------------------------------------
# First, grab a cell from the "regular" sheet that already carries the
# number/alignment styles we need (so we reuse existing style indices
# instead of minting new numFmt entries), then stamp those formats onto the
# "critical" sheet's B/D columns before writing values.
# ---------------------------------------------------------------------------
$regular.Cells.Item(2, 2).Copy() | Out-Null
$critical.Range("B2:B10").PasteSpecial(-4122) | Out-Null

$regular.Cells.Item(2, 4).Copy() | Out-Null
$critical.Range("D2:D10").PasteSpecial(-4122) | Out-Null

$critical.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Remove the extra rows (11-13) on the "critical" sheet - the target sheet
# only spans A1:E10.
# ---------------------------------------------------------------------------
$critical.Rows("11:13").Delete()

# ---------------------------------------------------------------------------
# Write the data rows onto both sheets.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $critical.Cells.Item($r, 1).Value = $data[0]
    $critical.Cells.Item($r, 2).Value = $data[1]
    $critical.Cells.Item($r, 3).Value = $data[2]
    $critical.Cells.Item($r, 4).Value = $data[3]
    $critical.Cells.Item($r, 5).Value = $data[4]

    $regular.Cells.Item($r, 1).Value = $data[0]
    $regular.Cells.Item($r, 2).Value = $data[1]
    $regular.Cells.Item($r, 3).Value = $data[2]
    $regular.Cells.Item($r, 4).Value = $data[3]
    $regular.Cells.Item($r, 5).Value = $data[4]
}

# ---------------------------------------------------------------------------
# Restore selection state to match the edited workbook (rows 3-7 were
# multi-selected on the "regular" sheet with the active cell on row 7).
# ---------------------------------------------------------------------------
$critical.Range("F9").Select() | Out-Null
$regular.Rows("3:7").Select() | Out-Null
